# Update AB de Villiers batting-innings stats (runs, balls, fours, sixes)
# Columns: A=playerName, B=teamName, C=runs, D=balls, E=fours, F=sixes
# The rows of match-by-match stats were reordered/updated per the latest
# activity pulled from the Excel form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the rows whose runs/balls/fours/sixes actually changed
# (rows 11 and 13 are untouched by this update).
$data = @{
    2  = @("56","43","5","0")
    3  = @("35","21","1","2")
    4  = @("39","36","4","0")
    5  = @("24","24","1","1")
    6  = @("15","12","1","1")
    7  = @("73","33","5","6")
    8  = @("2","5","0","0")
    9  = @("12","10","1","0")
    10 = @("0","2","0","0")
    12 = @("55","24","4","4")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 3).Value = $values[0]
    $ws.Cells.Item($row, 4).Value = $values[1]
    $ws.Cells.Item($row, 5).Value = $values[2]
    $ws.Cells.Item($row, 6).Value = $values[3]
}
